$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1): updates to column F ("想去人数") ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6995
$ws1.Range("F4").Value = 0
$ws1.Range("F5").Value = 457
$ws1.Range("F6").Value = 0
$ws1.Range("F7").Value = 0
$ws1.Range("F8").Value = 0
$ws1.Range("F11").Value = 0
$ws1.Range("F12").Value = 109
$ws1.Range("F16").Value = 416
$ws1.Range("F17").Value = 0
$ws1.Range("F18").Value = 0
$ws1.Range("F19").Value = 17
$ws1.Range("F20").Value = 5261
$ws1.Range("F21").Value = 121
$ws1.Range("F22").Value = 173
$ws1.Range("F23").Value = 0
$ws1.Range("F24").Value = 217
$ws1.Range("F25").Value = 243

# --- Sheet "演出" (sheet2): updates to column F ("想去人数") ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 0

# --- Sheet "全部类型" (sheet4): updates to column F ("想去人数") ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 0
$ws4.Range("F3").Value = 100
$ws4.Range("F4").Value = 0
$ws4.Range("F5").Value = 0
$ws4.Range("F6").Value = 0
$ws4.Range("F7").Value = 6895
$ws4.Range("F8").Value = 0
$ws4.Range("F9").Value = 202
$ws4.Range("F10").Value = 0
$ws4.Range("F11").Value = 0
$ws4.Range("F13").Value = 410
$ws4.Range("F15").Value = 18
$ws4.Range("F17").Value = 0
$ws4.Range("F18").Value = 0
$ws4.Range("F19").Value = 17
$ws4.Range("F21").Value = 0
$ws4.Range("F22").Value = 46
$ws4.Range("F23").Value = 121
$ws4.Range("F24").Value = 173
$ws4.Range("F25").Value = 672
$ws4.Range("F26").Value = 217
